$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.223.99"
$ws.Range("E2").Value = "  -5.12%  "
$ws.Range("D3").Value = "3.080.70"
$ws.Range("E3").Value = "  -5.45%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.085.46"
$ws.Range("E8").Value = "  -5.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("E10").Value = "  -6.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000215"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.04%  "
$ws.Range("D15").Value = "3.581.10"
$ws.Range("E15").Value = "  -5.41%  "
$ws.Range("D16").Value = "63.255.23"
$ws.Range("E16").Value = "  -5.12%  "
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "3.082.71"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.91%  "
$ws.Range("E28").Value = "  -8.30%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "59.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.94%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "499.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.77%  "
$ws.Range("E36").Value = "  -6.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.71%  "
$ws.Range("D38").Value = "3.147.27"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0394"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.19%  "
$ws.Range("E41").Value = "  -9.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -15.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.253"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("D50").Value = "0.0₃0503"
$ws.Range("E50").Value = "  -10.24%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.95%  "
